# Ran code for averaged intensities on spiral schemes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing rows 10-16 (column B) end up with new text because the
# averaging-scheme lookup table gained three new "Spiral" entries that
# were inserted right after "Gaussian-Quadrature".
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# New rows 17-19 for the remaining schemes that got pushed off the end.
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C17:P19").Value = 1

# Match the formatting of the existing index column (bold, centered,
# top-aligned, thin border) by copying the format from the row above.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
